$wb = $excel.ActiveWorkbook
$veiculos = $wb.Worksheets.Item(1)

# Create the new "Despesas" sheet positioned after "Veiculos"
$despesas = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $veiculos)
$despesas.Name = "Despesas"

# Populate "Despesas" header row + first data row first (matches author's
# original fill order so new shared strings land in the same slots)
$despesas.Cells.Item(1, 1).Value = "ID"
$despesas.Cells.Item(1, 2).Value = "ID Veículo"
$despesas.Cells.Item(1, 3).Value = "Tipo"
$despesas.Cells.Item(1, 4).Value = "Descrição"
$despesas.Cells.Item(1, 5).Value = "Data"
$despesas.Cells.Item(1, 6).Value = "Valor"

$despesas.Cells.Item(2, 1).Value = 1
$despesas.Cells.Item(2, 2).Value = 2
$despesas.Cells.Item(2, 3).Value = 1
$despesas.Cells.Item(2, 4).Value = "Falto bomba"
$despesas.Cells.Item(2, 5).Value = "22/05/2023"
$despesas.Cells.Item(2, 6).Value = 250

# Append the new vehicle row to "Veiculos"
$veiculos.Cells.Item(3, 1).Value = 6
$veiculos.Cells.Item(3, 2).Value = "QWE3D12"
$veiculos.Cells.Item(3, 3).Value = "Corvetezinho"
$veiculos.Cells.Item(3, 4).Value = "Chevrolet"
$veiculos.Cells.Item(3, 5).Value = 2025
$veiculos.Cells.Item(3, 6).Value = "Ativo"

# Finish populating "Despesas" with the second data row
$despesas.Cells.Item(3, 1).Value = 2
$despesas.Cells.Item(3, 2).Value = 6
$despesas.Cells.Item(3, 3).Value = 2
$despesas.Cells.Item(3, 4).Value = "Quebrou"
$despesas.Cells.Item(3, 5).Value = "25/02/2025"
$despesas.Cells.Item(3, 6).Value = 250

# Keep "Veiculos" as the active/selected sheet, as in the source workbook
$veiculos.Activate()
